$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Append the new data row (Data!A69:C69) ---
# Copy the format from the row above so the date cell keeps its existing
# date-formatted style (m/d/yyyy) instead of creating a brand new style.
$ws.Range("A68").Copy() | Out-Null
$ws.Range("A69").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A69").Value = 43967
$ws.Range("B69").Value = 8661
$ws.Range("C69").Value = 380

# Move the selection to the next empty row, as in the saved file, without
# leaving "Data" as the active/tabSelected sheet.
$ws.Range("A70").Select() | Out-Null

# --- Update the chart to cover the new row ---
$chartSheet = $wb.Worksheets.Item("Chart")
$co = $chartSheet.ChartObjects().Item(1)
$chart = $co.Chart

$sCases = $chart.SeriesCollection().Item(1)
$sCases.Formula = "=SERIES(Data!`$B`$1,Data!`$A`$2:`$A`$69,Data!`$B`$2:`$B`$69,1)"

$sDeaths = $chart.SeriesCollection().Item(2)
$sDeaths.Formula = "=SERIES(Data!`$C`$1,Data!`$A`$2:`$A`$69,Data!`$C`$2:`$C`$69,2)"

# Remove the data table under the plot area.
$chart.HasDataTable = $false

# Turn on minor gridlines for the value axis.
$valueAxis = $chart.Axes(2)
$valueAxis.HasMinorGridlines = $true

# Restore the Chart sheet as the active tab/selection.
$chartSheet.Activate() | Out-Null
